# Update dataset tracking spreadsheet: rename DS001 -> DS002 and append
# "_seed42" to the config/dataset filenames referencing this dataset.

$wb = $excel.ActiveWorkbook

# --- Sheet: Dataset Registry ---
$ws1 = $wb.Worksheets.Item("Dataset Registry")
$ws1.Range("A2").Value = "DS002"
$ws1.Range("B2").Value = "n10000_f_init5_cont0_disc5_sep6p6_seed42_config.yml"
$ws1.Range("C2").Value = "n10000_f_init5_cont0_disc5_sep6p6_seed42_dataset.csv"

# --- Sheet: Configuration Details ---
$ws2 = $wb.Worksheets.Item("Configuration Details")
$ws2.Range("A2").Value = "DS002"
$ws2.Range("B2").Value = "n10000_f_init5_cont0_disc5_sep6p6_seed42_config.yml"

# --- Sheet: Feature Separation Details ---
$ws3 = $wb.Worksheets.Item("Feature Separation Details")
$ws3.Range("A2").Value = "DS002"
$ws3.Range("A3").Value = "DS002"
$ws3.Range("A4").Value = "DS002"
$ws3.Range("A5").Value = "DS002"
$ws3.Range("A6").Value = "DS002"

# --- Sheet: File Metadata ---
$ws4 = $wb.Worksheets.Item("File Metadata")
$ws4.Range("A2").Value = "DS002"
$ws4.Range("B2").Value = "configs/data_generation/n10000_f_init5_cont0_disc5_sep6p6_seed42_config.yml"
$ws4.Range("C2").Value = "data/n10000_f_init5_cont0_disc5_sep6p6_seed42_dataset.csv"
